$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing text in place (same cells, new wording) ---
$ws.Range("B1").Value = "Test Case: Testing to see if createing kpis reaches data base with no errors"
$ws.Range("C2").Value = "Step 1: While logged out go to the create kpi page"
$ws.Range("D2").Value = "I am returned to the login page"
$ws.Range("C3").Value = "Step 2: Login as a service manager and go to the create kpi page"
$ws.Range("D3").Value = "I am redirected to the create KPI page"

# --- Fill in the new test steps / expected results for rows 4-6 ---
$ws.Range("C4").Value = "Step 3: Fill out a KPI for an employee on your teams"
$ws.Range("D4").Value = "A new KPI will be added to the database on that employee"
$ws.Range("C5").Value = "Step 4: Fill out a KPI for an employee that is not on your teams"
$ws.Range("D5").Value = "An error pops up saying that that access is denied"
$ws.Range("C6").Value = "Step 5: Leave Certain fields that are required blank"
$ws.Range("D6").Value = "An error pops up saying that some of the fields have been left blank"

# --- Update the view: active selection now sits on C5 ---
$ws.Range("C5").Select()
